$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 and 42: coin identities swap (TheSandbox <-> TrustWalletToken)
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.346"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.393.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.628.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3752"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3636"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08196"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.526"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.327"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.628.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06973"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.525"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.388.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.137"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.457"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.296"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.807.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.238"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.789"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.037"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02777"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2511"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08769"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.283"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.970"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  -2.32%  "
